$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: extend the table (rows 45-55) by copying the alternating row
# banding/format from the last two template rows (43=odd/shaded, 44=even),
# then filling in the new case data for each row. ---

$ws.Range("A43:AK43").Copy()
$ws.Range("A45:AK45").PasteSpecial(-4122)
$ws.Range("A44:AK44").Copy()
$ws.Range("A46:AK46").PasteSpecial(-4122)
$ws.Range("A43:AK43").Copy()
$ws.Range("A47:AK47").PasteSpecial(-4122)
$ws.Range("A44:AK44").Copy()
$ws.Range("A48:AK48").PasteSpecial(-4122)
$ws.Range("A43:AK43").Copy()
$ws.Range("A49:AK49").PasteSpecial(-4122)
$ws.Range("A44:AK44").Copy()
$ws.Range("A50:AK50").PasteSpecial(-4122)
$ws.Range("A43:AK43").Copy()
$ws.Range("A51:AK51").PasteSpecial(-4122)
$ws.Range("A44:AK44").Copy()
$ws.Range("A52:AK52").PasteSpecial(-4122)
$ws.Range("A43:AK43").Copy()
$ws.Range("A53:AK53").PasteSpecial(-4122)
$ws.Range("A44:AK44").Copy()
$ws.Range("A54:AK54").PasteSpecial(-4122)
$ws.Range("A43:AK43").Copy()
$ws.Range("A55:AK55").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 45
$ws.Range("A45").Value = 43
$ws.Range("B45").Value = "服務"
$ws.Range("C45").Value = 2025071382
$ws.Range("F45").Value = "L535"
$ws.Range("G45").Value = "田倉新樹店"
$ws.Range("H45").Value = "新北市新莊區"
$ws.Range("Q45").Value = "THILF0L535"
$ws.Range("R45").Value = "新北一"
$ws.Range("S45").Value = "湯家瑋"
$ws.Range("T45").Value = 1
$ws.Range("U45").Value = "已完工"
$ws.Range("V45").Value = "2025-07-10 11:20:02"
$ws.Range("W45").Value = "2025-07-10 11:00:00"
$ws.Range("X45").Value = "2025-07-10 11:20:00"
$ws.Range("Z45").Value = 0.3
$ws.Range("AB45").Value = "到場處理"
$ws.Range("AC45").Value = "PMQ3"
$ws.Range("AD45").Value = "O"
$ws.Range("AK45").Value = "O"

# Row 46
$ws.Range("A46").Value = 44
$ws.Range("B46").Value = "服務"
$ws.Range("C46").Value = 2025071393
$ws.Range("F46").Value = "L535"
$ws.Range("G46").Value = "田倉新樹店"
$ws.Range("H46").Value = "新北市新莊區"
$ws.Range("Q46").Value = "THILF0L535"
$ws.Range("R46").Value = "新北一"
$ws.Range("S46").Value = "湯家瑋"
$ws.Range("T46").Value = 1
$ws.Range("U46").Value = "已完工"
$ws.Range("V46").Value = "2025-07-10 11:30:40"
$ws.Range("W46").Value = "2025-07-10 11:00:00"
$ws.Range("X46").Value = "2025-07-10 11:20:00"
$ws.Range("Z46").Value = 0.3
$ws.Range("AB46").Value = "到場處理"
$ws.Range("AC46").Value = "PMQ3"
$ws.Range("AD46").Value = "O"
$ws.Range("AK46").Value = "O"

# Row 47
$ws.Range("A47").Value = 45
$ws.Range("B47").Value = "服務"
$ws.Range("C47").Value = 2025071404
$ws.Range("F47").Value = 2109
$ws.Range("G47").Value = "新莊瓊林南"
$ws.Range("H47").Value = "新北市新莊區"
$ws.Range("Q47").Value = "THILF02109"
$ws.Range("R47").Value = "新北一"
$ws.Range("S47").Value = "湯家瑋"
$ws.Range("T47").Value = 1
$ws.Range("U47").Value = "已完工"
$ws.Range("V47").Value = "2025-07-10 11:42:58"
$ws.Range("W47").Value = "2025-07-10 11:25:00"
$ws.Range("X47").Value = "2025-07-10 11:45:00"
$ws.Range("Z47").Value = 0.3
$ws.Range("AB47").Value = "到場處理"
$ws.Range("AC47").Value = "PMQ3"
$ws.Range("AD47").Value = "O"
$ws.Range("AK47").Value = "O"

# Row 48
$ws.Range("A48").Value = 46
$ws.Range("B48").Value = "服務"
$ws.Range("C48").Value = 2025071415
$ws.Range("F48").Value = "D070"
$ws.Range("G48").Value = "三重重新店"
$ws.Range("H48").Value = "新北市三重區"
$ws.Range("Q48").Value = "THILF0D070"
$ws.Range("R48").Value = "新北一"
$ws.Range("S48").Value = "吳宗鴻"
$ws.Range("T48").Value = 1
$ws.Range("U48").Value = "已完工"
$ws.Range("V48").Value = "2025-07-10 12:05:57"
$ws.Range("W48").Value = "2025-07-10 11:50:00"
$ws.Range("X48").Value = "2025-07-10 12:05:00"
$ws.Range("Z48").Value = 0.3
$ws.Range("AB48").Value = "到場處理"
$ws.Range("AC48").Value = "PMQ3"
$ws.Range("AD48").Value = "O"
$ws.Range("AK48").Value = "O"

# Row 49
$ws.Range("A49").Value = 47
$ws.Range("B49").Value = "服務"
$ws.Range("C49").Value = 2025071417
$ws.Range("F49").Value = 3627
$ws.Range("G49").Value = "三重重新橋"
$ws.Range("H49").Value = "新北市三重區"
$ws.Range("Q49").Value = "THILF03627"
$ws.Range("R49").Value = "新北一"
$ws.Range("S49").Value = "吳宗鴻"
$ws.Range("T49").Value = 1
$ws.Range("U49").Value = "已完工"
$ws.Range("V49").Value = "2025-07-10 12:27:26"
$ws.Range("W49").Value = "2025-07-10 12:12:00"
$ws.Range("X49").Value = "2025-07-10 12:27:00"
$ws.Range("Z49").Value = 0.3
$ws.Range("AB49").Value = "到場處理"
$ws.Range("AC49").Value = "PMQ3"
$ws.Range("AD49").Value = "O"
$ws.Range("AK49").Value = "O"

# Row 50
$ws.Range("A50").Value = 48
$ws.Range("B50").Value = "服務"
$ws.Range("C50").Value = 2025071418
$ws.Range("F50").Value = 3627
$ws.Range("G50").Value = "三重重新橋"
$ws.Range("H50").Value = "新北市三重區"
$ws.Range("Q50").Value = "THILF03627"
$ws.Range("R50").Value = "新北一"
$ws.Range("S50").Value = "吳宗鴻"
$ws.Range("T50").Value = 1
$ws.Range("U50").Value = "已完工"
$ws.Range("V50").Value = "2025-07-10 12:29:47"
$ws.Range("W50").Value = "2025-07-10 12:10:00"
$ws.Range("X50").Value = "2025-07-10 12:25:00"
$ws.Range("Z50").Value = 0.3
$ws.Range("AB50").Value = "到場處理"
$ws.Range("AC50").Value = "TVV"
$ws.Range("AD50").Value = "O"
$ws.Range("AK50").Value = "O"

# Row 51
$ws.Range("A51").Value = 49
$ws.Range("B51").Value = "服務"
$ws.Range("C51").Value = 2025071425
$ws.Range("F51").Value = 2259
$ws.Range("G51").Value = "三重興華店"
$ws.Range("H51").Value = "新北市三重區"
$ws.Range("Q51").Value = "THILF02259"
$ws.Range("R51").Value = "新北一"
$ws.Range("S51").Value = "吳宗鴻"
$ws.Range("T51").Value = 1
$ws.Range("U51").Value = "已完工"
$ws.Range("V51").Value = "2025-07-10 13:15:23"
$ws.Range("W51").Value = "2025-07-10 12:50:00"
$ws.Range("X51").Value = "2025-07-10 13:05:00"
$ws.Range("Z51").Value = 0.3
$ws.Range("AB51").Value = "到場處理"
$ws.Range("AC51").Value = "PMQ3"
$ws.Range("AD51").Value = "O"
$ws.Range("AK51").Value = "O"

# Row 52
$ws.Range("A52").Value = 50
$ws.Range("B52").Value = "服務"
$ws.Range("C52").Value = 2025071430
$ws.Range("F52").Value = 5291
$ws.Range("G52").Value = "新莊國家置地"
$ws.Range("H52").Value = "新北市新莊區"
$ws.Range("Q52").Value = "THILF05291"
$ws.Range("R52").Value = "新北一"
$ws.Range("S52").Value = "湯家瑋"
$ws.Range("T52").Value = 1
$ws.Range("U52").Value = "已完工"
$ws.Range("V52").Value = "2025-07-10 13:33:23"
$ws.Range("W52").Value = "2025-07-10 13:00:00"
$ws.Range("X52").Value = "2025-07-10 13:30:00"
$ws.Range("Z52").Value = 0.5
$ws.Range("AB52").Value = "到場處理"
$ws.Range("AC52").Value = "PMQ3+TVV"
$ws.Range("AD52").Value = "O"
$ws.Range("AJ52").Value = "O"
$ws.Range("AK52").Value = "O"

# Row 53
$ws.Range("A53").Value = 51
$ws.Range("B53").Value = "服務"
$ws.Range("C53").Value = 2025071446
$ws.Range("F53").Value = 2837
$ws.Range("G53").Value = "北縣重萬店"
$ws.Range("H53").Value = "新北市三重區"
$ws.Range("Q53").Value = "THILF02837"
$ws.Range("R53").Value = "新北一"
$ws.Range("S53").Value = "吳宗鴻"
$ws.Range("T53").Value = 1
$ws.Range("U53").Value = "已完工"
$ws.Range("V53").Value = "2025-07-10 14:35:49"
$ws.Range("W53").Value = "2025-07-10 14:20:00"
$ws.Range("X53").Value = "2025-07-10 14:35:00"
$ws.Range("Z53").Value = 0.3
$ws.Range("AB53").Value = "到場處理"
$ws.Range("AC53").Value = "PMQ3"
$ws.Range("AD53").Value = "O"
$ws.Range("AK53").Value = "O"

# Row 54
$ws.Range("A54").Value = 52
$ws.Range("B54").Value = "服務"
$ws.Range("C54").Value = 2025071450
$ws.Range("F54").Value = 3676
$ws.Range("G54").Value = "三重田心店"
$ws.Range("H54").Value = "新北市三重區"
$ws.Range("Q54").Value = "THILF03676"
$ws.Range("R54").Value = "新北一"
$ws.Range("S54").Value = "吳宗鴻"
$ws.Range("T54").Value = 1
$ws.Range("U54").Value = "已完工"
$ws.Range("V54").Value = "2025-07-10 15:00:25"
$ws.Range("W54").Value = "2025-07-10 14:45:00"
$ws.Range("X54").Value = "2025-07-10 15:00:00"
$ws.Range("Z54").Value = 0.3
$ws.Range("AB54").Value = "到場處理"
$ws.Range("AC54").Value = "PMQ3"
$ws.Range("AD54").Value = "O"
$ws.Range("AK54").Value = "O"

# Row 55
$ws.Range("A55").Value = 53
$ws.Range("B55").Value = "服務"
$ws.Range("C55").Value = 2025071465
$ws.Range("F55").Value = "D024"
$ws.Range("G55").Value = "三重三民店"
$ws.Range("H55").Value = "新北市三重區"
$ws.Range("Q55").Value = "THILF0D024"
$ws.Range("R55").Value = "新北一"
$ws.Range("S55").Value = "吳宗鴻"
$ws.Range("T55").Value = 1
$ws.Range("U55").Value = "已完工"
$ws.Range("V55").Value = "2025-07-10 15:30:32"
$ws.Range("W55").Value = "2025-07-10 15:14:00"
$ws.Range("X55").Value = "2025-07-10 15:29:00"
$ws.Range("Z55").Value = 0.3
$ws.Range("AB55").Value = "到場處理"
$ws.Range("AC55").Value = "PMQ3"
$ws.Range("AD55").Value = "O"
$ws.Range("AK55").Value = "O"

# --- Step 2: the source data also normalizes P44/AC44 to the wrap-text
# style used elsewhere in the even-row band (cosmetic parity with the
# newly appended rows). ---
$ws.Range("P44").WrapText = $true
$ws.Range("AC44").WrapText = $true

# --- Step 3: grow the print area to cover the newly added rows. ---
$printAreaName = $wb.Names.Item(1)
$printAreaName.RefersTo = "='Report'!`$A`$1:`$AK`$55"

# --- Step 4: move the active selection to the new last row, matching
# where the author left off after data entry. ---
$ws.Range("A55").Select()
